$wb = $excel.ActiveWorkbook

# --- Sheet "nokey1": insert a "# of Diffs" column before the existing data ---
$ws1 = $wb.Worksheets.Item("nokey1")
$ws1.AutoFilterMode = $false

$ws1.Columns.Item(1).Insert()

$ws1.Range("B1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)
$ws1.Range("A1").Value = "# of Diffs"

$ws1.Range("B2").Copy()
$ws1.Range("A2").PasteSpecial(-4122)
$ws1.Range("A2").Value = 1

$excel.CutCopyMode = 0

$ws1.Columns.Item(1).ColumnWidth = 13.15
$ws1.Columns.Item(2).ColumnWidth = 15.32
$ws1.Columns.Item(3).ColumnWidth = 15.32

$null = $ws1.Range("A1:C2").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "nokey1!_FilterDatabase") {
        $n.RefersTo = "=nokey1!`$A`$1:`$C`$2"
    }
}

# --- Sheet "nokey2": insert a "# of Diffs" column before the existing data ---
$ws2 = $wb.Worksheets.Item("nokey2")
$ws2.AutoFilterMode = $false

$ws2.Columns.Item(1).Insert()

$ws2.Range("B1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A1").Value = "# of Diffs"

$ws2.Range("A2").Value = 0

$excel.CutCopyMode = 0

$ws2.Columns.Item(1).ColumnWidth = 13.15
$ws2.Columns.Item(2).ColumnWidth = 15.32
$ws2.Columns.Item(3).ColumnWidth = 15.32

$null = $ws2.Range("A1:C2").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "nokey2!_FilterDatabase") {
        $n.RefersTo = "=nokey2!`$A`$1:`$C`$2"
    }
}
